{"js": "// The commit fixes the casing of the placeholder token used in the\n// \"file naming format\" sentence: \"WORLDNAME\" -> \"Worldname\"\n// (the run keeps its existing italic formatting; only the text changes).\n\nconst searchResults = context.document.body.search(\"WORLDNAME\", {\n    matchCase: true,\n    matchWholeWord: true\n});\nsearchResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n    searchResults.items[i].insertText(\"Worldname\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The commit fixes the casing of the placeholder token used in the\n# \"file naming format\" sentence: \"WORLDNAME\" -> \"Worldname\"\n# (the run keeps its existing italic formatting; only the text changes).\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"WORLDNAME\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $true\n$find.Replacement.Text = \"Worldname\"\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll)\n"}
